$wb = $excel.ActiveWorkbook

# Rename the "SwateTemplateMetadata" sheet to "isa_template"
$sheet = $wb.Worksheets.Item("SwateTemplateMetadata")
$sheet.Name = "isa_template"

# Update the active selection on that sheet from E11 to B14
$sheet.Activate()
$sheet.Range("B14").Select()
